# Lab5/Tests.xlsx - "Tuned catapult with solidified hardware design"
# Adds 8 new catapult test runs (rows 8-15) to the data already present in
# columns I:O (Launch_Angle / Acceleration / Range trial columns), and
# updates the sheet selection to cover the new data range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: 75, 2850, 5.5, 5.5, 5.5, 5.5, 5.5
$ws.Range("I8").Value = 75
$ws.Range("J8").Value = 2850
$ws.Range("K8").Value = 5.5
$ws.Range("L8").Value = 5.5
$ws.Range("M8").Value = 5.5
$ws.Range("N8").Value = 5.5
$ws.Range("O8").Value = 5.5

# Row 9: 65, 2750, 5.25, 5.25, 5.25, 5.25, 5.25
$ws.Range("I9").Value = 65
$ws.Range("J9").Value = 2750
$ws.Range("K9").Value = 5.25
$ws.Range("L9").Value = 5.25
$ws.Range("M9").Value = 5.25
$ws.Range("N9").Value = 5.25
$ws.Range("O9").Value = 5.25

# Row 10: 50, 2750, 5, 5, 5.25, 5, 5.25
$ws.Range("I10").Value = 50
$ws.Range("J10").Value = 2750
$ws.Range("K10").Value = 5
$ws.Range("L10").Value = 5
$ws.Range("M10").Value = 5.25
$ws.Range("N10").Value = 5
$ws.Range("O10").Value = 5.25

# Row 11: 50, 2500, 5, 5, 5, 5, 5
$ws.Range("I11").Value = 50
$ws.Range("J11").Value = 2500
$ws.Range("K11").Value = 5
$ws.Range("L11").Value = 5
$ws.Range("M11").Value = 5
$ws.Range("N11").Value = 5
$ws.Range("O11").Value = 5

# Row 12: 40, 2500, 5.25, 5.25, 5.25, 5, 5
$ws.Range("I12").Value = 40
$ws.Range("J12").Value = 2500
$ws.Range("K12").Value = 5.25
$ws.Range("L12").Value = 5.25
$ws.Range("M12").Value = 5.25
$ws.Range("N12").Value = 5
$ws.Range("O12").Value = 5

# Row 13: 40, 2250, 2.5, 2.5, 3, 2, 2.5
$ws.Range("I13").Value = 40
$ws.Range("J13").Value = 2250
$ws.Range("K13").Value = 2.5
$ws.Range("L13").Value = 2.5
$ws.Range("M13").Value = 3
$ws.Range("N13").Value = 2
$ws.Range("O13").Value = 2.5

# Row 14: 40, 2400, 4.75, 3, 4.5, 4.5, 4.5
$ws.Range("I14").Value = 40
$ws.Range("J14").Value = 2400
$ws.Range("K14").Value = 4.75
$ws.Range("L14").Value = 3
$ws.Range("M14").Value = 4.5
$ws.Range("N14").Value = 4.5
$ws.Range("O14").Value = 4.5

# Row 15: 40, 2475, 4.75, 4.75, 4.75, 4.75, 4.75
$ws.Range("I15").Value = 40
$ws.Range("J15").Value = 2475
$ws.Range("K15").Value = 4.75
$ws.Range("L15").Value = 4.75
$ws.Range("M15").Value = 4.75
$ws.Range("N15").Value = 4.75
$ws.Range("O15").Value = 4.75

# Select the full data block (I1:O15) to match the final selection state
# recorded in the workbook after the edit.
$ws.Range("I1:O15").Select()
